$wb = $excel.ActiveWorkbook

# --- Sheet "2025" (sheet1) row 2 updates ---
$ws2025 = $wb.Worksheets.Item("2025")
$ws2025.Range("B2").Value = 5457.008988199987
$ws2025.Range("E2").Value = 212716.0964344695
$ws2025.Range("I2").Value = 109810.98647904
$ws2025.Range("L2").Value = 428396.954927634
$ws2025.Range("M2").Value = 80191.32604175001
$ws2025.Range("N2").Value = 44753.37472369202
$ws2025.Range("O2").Value = 50658.98282837399

# --- Sheet "2030" (sheet2) row 2 updates ---
$ws2030 = $wb.Worksheets.Item("2030")
$ws2030.Range("A2").Value = 8830.954041229934
$ws2030.Range("B2").Value = 40368.76755562462
$ws2030.Range("E2").Value = 288326.2489419255
$ws2030.Range("I2").Value = 264868.335261012
$ws2030.Range("L2").Value = 217212.4395345119
$ws2030.Range("M2").Value = 123469.0689506697
$ws2030.Range("N2").Value = 63706.06971442258
$ws2030.Range("O2").Value = 54067.35526560284
